$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step A: shift existing J/K columns right into K/M, reusing the
#     existing shared strings (J -> K, old K -> M). Capture old K values
#     first since we are about to overwrite column K with column J's data.
$oldK1 = $ws.Range("K1").Value2
$oldK2 = $ws.Range("K2").Value2
$oldK3 = $ws.Range("K3").Value2
$oldK4 = $ws.Range("K4").Value2
$oldK5 = $ws.Range("K5").Value2

$oldJ1 = $ws.Range("J1").Value2
$oldJ2 = $ws.Range("J2").Value2
$oldJ3 = $ws.Range("J3").Value2
$oldJ4 = $ws.Range("J4").Value2
$oldJ5 = $ws.Range("J5").Value2

$ws.Range("M1").Value = $oldK1
$ws.Range("M2").Value = $oldK2
$ws.Range("M3").Value = $oldK3
$ws.Range("M4").Value = $oldK4
$ws.Range("M5").Value = $oldK5

$ws.Range("K1").Value = $oldJ1
$ws.Range("K2").Value = $oldJ2
$ws.Range("K3").Value = $oldJ3
$ws.Range("K4").Value = $oldJ4
$ws.Range("K5").Value = $oldJ5

$ws.Range("J2:J4").Clear()

# --- Step B: new "States" column (L), top to bottom.
$ws.Range("L1").Value = "States"
$ws.Range("L2").Value = "AK"
$ws.Range("L3").Value = "ON"
$ws.Range("L5").Value = "Kerala"

# --- Step C: new "Number Of Learners" column (J).
$ws.Range("J1").Value = "Number Of Learners"
$ws.Range("J5").Value = "01-04"

# --- Formatting ---
# D2:D5 phone numbers become Text-formatted.
$ws.Range("D2:D5").NumberFormat = "@"
# The "01-04" entry picks up the same Text format (prevents date coercion).
$ws.Range("J5").NumberFormat = "@"

# Header row 1 grows taller to fit the new wrapped column headers.
$ws.Rows("1").RowHeight = 61.8

# Selection moves to D2.
$ws.Range("D2").Select() | Out-Null
